# diamonds is_deleted from excel sheet
# - remove the now-unused "Cut Grade" header / row of placeholder formatting
# - rename/re-purpose trailing "image-*" headers into the new
#   COMMENT / LOCATION / image-1 / is_deleted columns
# - refresh the header row styling (Calibri 10, centered, thin borders)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 only ever held stray cell formatting (no values) - drop it entirely.
$ws.Rows("2:2").Delete()

# "Cut Grade" -> "Cut"
$ws.Range("I1").Value = "Cut"

# Shift the trailing placeholder image-* columns into their new roles.
$ws.Range("Z1").Value = "COMMENT"
$ws.Range("AA1").Value = "LOCATION"
$ws.Range("AB1").Value = "image-1"
$ws.Range("AC1").Value = "is_deleted"

# Refresh header styling: Calibri 10pt (no bold, automatic/theme text
# color), centered, thin black borders. Set the theme color *before* the
# other font properties so every header cell converges onto one shared
# font entry (some header cells previously had an explicit black font
# color, others inherited the theme color).
$header = $ws.Range("A1:AB1")
$header.Font.ThemeColor = 1
$header.Font.Name = "Calibri"
$header.Font.Size = 10
$header.Font.Bold = $false
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4107
$header.Borders.Item(7).LineStyle = 1
$header.Borders.Item(7).Weight = 2
$header.Borders.Item(7).Color = 0
$header.Borders.Item(8).LineStyle = 1
$header.Borders.Item(8).Weight = 2
$header.Borders.Item(8).Color = 0
$header.Borders.Item(9).LineStyle = 1
$header.Borders.Item(9).Weight = 2
$header.Borders.Item(9).Color = 0
$header.Borders.Item(10).LineStyle = 1
$header.Borders.Item(10).Weight = 2
$header.Borders.Item(10).Color = 0

# Last header cell (is_deleted) only keeps left/right borders.
$lastCell = $ws.Range("AC1")
$lastCell.Font.ThemeColor = 1
$lastCell.Font.Name = "Calibri"
$lastCell.Font.Size = 10
$lastCell.Font.Bold = $false
$lastCell.HorizontalAlignment = -4108
$lastCell.VerticalAlignment = -4107
$lastCell.Borders.Item(7).LineStyle = 1
$lastCell.Borders.Item(7).Weight = 2
$lastCell.Borders.Item(7).Color = 0
$lastCell.Borders.Item(10).LineStyle = 1
$lastCell.Borders.Item(10).Weight = 2
$lastCell.Borders.Item(10).Color = 0
$lastCell.Borders.Item(8).LineStyle = -4142
$lastCell.Borders.Item(9).LineStyle = -4142

# Park the selection back at A1 (the sheet no longer scrolls to column I).
$ws.Range("A1").Select()

$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 1
